# Apply the LinuxForHealth re-brand / regeneration updates to the
# StructureDefinition-reference-code workbook.
#
# Sheet "Metadata": update URL, Version, Date and Publisher values.
# Sheet "Elements": the base "Extension" row's Constraint(s) cell (AI2) is
# cleared - the ele-1/ext-1 constraint text now lives only on the
# "Extension.extension" row (AI4), which already carries it and needs no
# change.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/reference-code"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AI2").Value = ""
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/reference-code"
